# Insert a new data row before the current row 143 (Feria Lagunitas de
# Puerto Montt - Brócoli weekly price records). This shifts the existing
# rows 143..276 down to 144..277, matching the target diff which grows the
# sheet from A1:R276 to A1:R277 and re-numbers every subsequent record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(143).Insert()

# Populate the freshly inserted row 143 with the new record.
$row = 143
$ws.Cells.Item($row, 1).Value  = 4
$ws.Cells.Item($row, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value  = "Los Lagos"
$ws.Cells.Item($row, 4).Value  = 44589
$ws.Cells.Item($row, 5).Value  = 10
$ws.Cells.Item($row, 6).Value  = 100112023
$ws.Cells.Item($row, 7).Value  = "Brócoli"
$ws.Cells.Item($row, 8).Value  = "Sin especificar"
$ws.Cells.Item($row, 9).Value  = "Primera"
$ws.Cells.Item($row, 10).Value = 500
$ws.Cells.Item($row, 11).Value = 1500
$ws.Cells.Item($row, 12).Value = 1500
$ws.Cells.Item($row, 13).Value = 1500
$ws.Cells.Item($row, 14).Value = "$/unidad"
$ws.Cells.Item($row, 15).Value = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value = 1500
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
